{"js": "// \"Private fields, with public getter and setter functions for more secure\n// data control\" -> \"Protected fields, with public getter and setter\n// functions for more secure data control\"\n// (i.e. \"Private\" becomes \"Protected\"; the rest of the sentence is unchanged)\n\nconst body = context.document.body;\n\nconst results = body.search(\n  \"Private fields, with public getter and setter functions for more secure data control\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not locate the target paragraph text to edit.\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\n  \"Protected fields, with public getter and setter functions for more secure data control\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# \"Private fields, with public getter and setter functions for more secure\n# data control\" -> \"Protected fields, with public getter and setter\n# functions for more secure data control\"\n# (i.e. \"Private\" becomes \"Protected\"; the rest of the sentence is unchanged)\n\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Private fields, with public getter and setter functions for more secure data control\"\n$find.Replacement.Text = \"Protected fields, with public getter and setter functions for more secure data control\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
